$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scheduled cryptocurrency price/volume refresh (GitHub Actions bot).
# Column D holds price text that sometimes *looks* numeric (e.g. "215.82");
# prefixing with a literal apostrophe keeps Excel from coercing those
# cells into real numbers, matching the source data which stores them as text.

$ws.Range("D2").Value = "27.020.36"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "1.678.85"
$ws.Range("E3").Value = "  +0.82%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'215.82"
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("E6").Value = "  -2.43%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "'0.254"
$ws.Range("E8").Value = "  +2.02%  "
$ws.Range("D9").Value = "'21.39"
$ws.Range("E9").Value = "  +5.44%  "
$ws.Range("E10").Value = "  +0.64%  "
$ws.Range("E11").Value = "  -0.96%  "
$ws.Range("D12").Value = "1.916.07"
$ws.Range("E12").Value = "  +0.85%  "
$ws.Range("D13").Value = "1.726.46"
$ws.Range("E13").Value = "  +3.85%  "
$ws.Range("E14").Value = "  +0.80%  "
$ws.Range("E15").Value = "  +1.72%  "
$ws.Range("D16").Value = "'66.49"
$ws.Range("E16").Value = "  +0.47%  "
$ws.Range("D17").Value = "27.017.64"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").Value = "'8.16"
$ws.Range("E18").Value = "  +1.96%  "
$ws.Range("D19").Value = "'235.85"
$ws.Range("E19").Value = "  +0.83%  "
$ws.Range("D20").Value = "0.0₃0737"
$ws.Range("E20").Value = "  +0.72%  "
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("E22").Value = "  +2.15%  "
$ws.Range("E23").Value = "  +1.57%  "
$ws.Range("D24").Value = "'2.12"
$ws.Range("E24").Value = "  -4.27%  "
$ws.Range("D25").Value = "'146.48"
$ws.Range("E26").Value = "  +2.00%  "
$ws.Range("D27").Value = "'16.41"
$ws.Range("E27").Value = "  +3.21%  "
$ws.Range("E28").Value = "  -2.39%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("E30").Value = "  +0.20%  "
$ws.Range("E31").Value = "  +0.17%  "
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("D33").Value = "1.536.48"
$ws.Range("E34").Value = "  +1.44%  "
$ws.Range("E35").Value = "  +5.26%  "
$ws.Range("E36").Value = "  -0.76%  "
$ws.Range("D37").Value = "'0.589"
$ws.Range("E37").Value = "  +1.72%  "
$ws.Range("D38").Value = "'0.919"
$ws.Range("E38").Value = "  +1.61%  "
$ws.Range("E39").Value = "  +3.22%  "
$ws.Range("E40").Value = "  +6.71%  "
$ws.Range("D42").Value = "'67.99"
$ws.Range("E42").Value = "  +3.14%  "
$ws.Range("D43").Value = "'5.61"
$ws.Range("E43").Value = "  -2.36%  "
$ws.Range("E44").Value = "  -0.31%  "
$ws.Range("D45").Value = "1.819.44"
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("E46").Value = "  -0.37%  "
$ws.Range("D47").Value = "'90.41"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("E49").Value = "  +2.51%  "
$ws.Range("D50").Value = "'7.96"
$ws.Range("E50").Value = "  +5.44%  "
$ws.Range("D51").Value = "'0.0505"
$ws.Range("E51").Value = "  -0.08%  "
